$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold, centered, bordered) from the last
# existing header cell (AC1) onto the three new header cells before
# writing their text, so they match the rest of the header row style.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# New header labels for the season-record columns.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Populate every data row (2-50) with the team's season record.
$ws.Range("AD2:AD50").Value = 73
$ws.Range("AE2:AE50").Value = 89
$ws.Range("AF2:AF50").Value = 0

Write-Host "Added Wins/Losses/Ties columns (AD:AF) for rows 1-50"
